# Refactor synthetic "status" array: swap square-emoji glyphs for
# book-emoji glyphs, and rename the "noir" (black) status label to
# "bleu" (blue) to match the new blue-book glyph.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old glyph -> new glyph (column A, "statut")
$glyphMap = @{
    "🟥" = "📕"
    "🟧" = "📙"
    "⬛" = "📘"
    "🟩" = "📗"
}

# Mapping of old label -> new label (column B, "statut_label")
$labelMap = @{
    "noir" = "bleu"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value2
    if ($glyphMap.ContainsKey($valA)) {
        $cellA.Value2 = $glyphMap[$valA]
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value2
    if ($labelMap.ContainsKey($valB)) {
        $cellB.Value2 = $labelMap[$valB]
    }
}
